$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Defined names (workbook-level) ---
$wb.Names.Add('deeznuts', '=OFFSET(Sheet1!$A$1,0,0,COUNTA(Sheet1!$A:$A),1)')
$wb.Names.Add('URMUM', '=OFFSET(Sheet1!$G$8,0,0,COUNTA(Sheet1!$G:$G),1)')
$wb.Names.Add('urnutz', '=OFFSET(Sheet1!$D$1,0,0,COUNTA(Sheet1!$D:$D),1)')

# --- Column D (1..13 running count) ---
$ws.Range("D1").Value = 3
$ws.Range("D2").Value = 4
$ws.Range("D3").Value = 5
$ws.Range("D4").Value = 6
$ws.Range("D5").Value = 7
$ws.Range("D6").Value = 8
$ws.Range("D7").Value = 9
$ws.Range("D8").Value = 10
$ws.Range("D9").Value = 11
$ws.Range("D10").Value = 12
$ws.Range("D11").Value = 13
$ws.Range("D12").Value = 13
$ws.Range("D13").Value = 23

# --- Column A extra rows ---
$ws.Range("A11").Value = 13
$ws.Range("A12").Value = 13
$ws.Range("A13").Value = 23

# --- B2 formula rewritten to inline the B1 formula ---
$ws.Range("B2").Formula = "=EOMONTH(DATE(2015,2,2),1)"

# --- Column G (room-nights style series, rows 8-24) ---
$ws.Range("G8").Value = 3
$ws.Range("G9").Value = 4
$ws.Range("G10").Value = 5
$ws.Range("G11").Value = 6
$ws.Range("G12").Value = 7
$ws.Range("G13").Value = 8
$ws.Range("G14").Value = 9
$ws.Range("G15").Value = 10
$ws.Range("G16").Value = 11
$ws.Range("G17").Value = 12
$ws.Range("G18").Value = 13
$ws.Range("G19").Value = 13
$ws.Range("G20").Value = 23
$ws.Range("G21").Value = 42
$ws.Range("G22").Value = 23
$ws.Range("G23").Value = 1
$ws.Range("G24").Value = 10

# --- Formulas that reference the new named ranges ---
$ws.Range("B8").Formula = "=SUM(urnutz)"
$ws.Range("I8").Formula = "=SUM(URMUM)"

# --- Selection moves to G25 ---
$ws.Range("G25").Select() | Out-Null
